$d = $word.ActiveDocument
$d.Content.Find.Execute("{{BACKING_AMOUNT_TEXT}}", $true, $false, $false, $false, $false,
                         $true, 1, $false, "{{BACKING_AMOUNT_TEXT}}", 2)
